# Generate Report for Handback
# Refresh the timestamp columns on each sheet to reflect a newer handback
# report generation run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 354a4a40-...-fefad.md row (row 3).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-10-13 13:08:17"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 354a4a40-...-fefad row (row 3).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-10-13 13:08:07"
$wsZhCn.Range("K3").Value = "2016-10-13 13:08:51"

# de-de sheet: same pair of columns for its row 3.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-10-13 13:08:17"
$wsDeDe.Range("K3").Value = "2016-10-13 13:09:08"
